# This script reproduces the edit described by the diff:
# a new data row was inserted at spreadsheet row 90 (pushing the
# previously-existing rows 90..184 down to 91..185), and the new row 90
# was populated with a fresh "Zapallo italiano" price record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 90; Excel shifts rows 90:184 down to 91:185
# and extends the used range / dimension to A1:R185 automatically.
$ws.Rows.Item(90).Insert()

# Populate the newly inserted row 90 with the new record's data.
$ws.Range("A90").Value = 7
$ws.Range("B90").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C90").Value = "Ñuble"
$ws.Range("D90").Value = 44586
$ws.Range("E90").Value = 16
$ws.Range("F90").Value = 100112032
$ws.Range("G90").Value = "Zapallo italiano"
$ws.Range("H90").Value = "Sin especificar"
$ws.Range("I90").Value = "Primera"
$ws.Range("J90").Value = 120
$ws.Range("K90").Value = 8500
$ws.Range("L90").Value = 9000
$ws.Range("M90").Value = 8750
$ws.Range("N90").Value = "`$/caja 60 unidades"
$ws.Range("O90").Value = "Región del Maule"
$ws.Range("P90").Value = 146
$ws.Range("Q90").Value = 60
$ws.Range("R90").Value = "Hortaliza"
